$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top of the data (row 2) for the newest date,
# pushing all existing rows down by one (Excel duplicates the row
# immediately below into the new bottom row, extending the table by one).
$ws.Rows.Item(2).Insert()

# Populate the new row with the latest price data. The date is entered
# with a leading apostrophe so it is stored as literal text (matching
# the rest of the date column) instead of being auto-converted into a
# date serial number.
$ws.Range("A2").Value = "'2026-01-30"
$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610
